# Adjust database input files to new standards:
#  - Add descriptions everywhere (row 5) and format/scenario notes (row 6)
#  - Replace the [p.u.] unit with [%, 0-1] for FirmCapCoef (column L)
#  - Widen YearCom/YearDecom columns (P:Q) to fit the new descriptions
#  - Grow row 5's height to fit the longer wrapped description text
#  - Reword the "Data Package" description

$wb = $excel.ActiveWorkbook

$descriptions = @{
    "C" = "Name of generator"
    "D" = "Corresponding technology"
    "E" = "Node where generator is connected"
    "F" = "Number of existing VRE units"
    "G" = "Maximum active power output of VRE unit"
    "H" = "Whether the model can invest in additional units (1) or not (0)"
    "I" = "Maximum number of VRE units that can be invested in"
    "J" = "Annualized investment cost per MW"
    "K" = "Operation and maintenance cost of VRE unit"
    "L" = "Firm capacity coefficient of VRE unit"
    "M" = "Maximum reactive power output of VRE unit"
    "N" = "Minimum reactive power output of VRE unit"
    "O" = "Inertia constant H of VRE unit"
    "P" = "Year where it is commissioned (1.1.xxxx)"
    "Q" = "Year where it is decommissioned (31.12.xxxx)"
    "R" = "Latitude"
    "S" = "Longitude"
    "T" = "Which package this belongs to"
    "U" = "Where the data for the entry comes from"
}

$formats = @{
    "C" = "-"
    "D" = "-"
    "E" = "-"
    "F" = "Scenario-dependent"
    "G" = "Scenario-dependent"
    "H" = "Scenario-dependent"
    "I" = "Scenario-dependent"
    "J" = "Scenario-dependent"
    "K" = "Scenario-dependent"
    "L" = "Scenario-dependent"
    "M" = "Scenario-dependent"
    "N" = "Scenario-dependent"
    "O" = "Scenario-dependent"
    "P" = "Scenario-dependent"
    "Q" = "Scenario-dependent"
    "R" = "-"
    "S" = "-"
    "T" = "Scenario-dependent"
    "U" = "Scenario-dependent"
}

foreach ($sheetName in @("ScenarioA", "ScenarioB")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 5: per-column descriptions
    foreach ($col in $descriptions.Keys) {
        $ws.Range($col + "5").Value = $descriptions[$col]
    }

    # Row 6: per-column format / scenario-dependency notes
    foreach ($col in $formats.Keys) {
        $ws.Range($col + "6").Value = $formats[$col]
    }

    # Row 7: FirmCapCoef is a share/percentage, not a p.u. quantity -> [%, 0-1]
    $ws.Range("L7").Value = "[%, 0-1]"

    # Row 5 needs to be taller to fit the new wrapped description text
    $ws.Rows.Item(5).RowHeight = 90

    # YearCom/YearDecom columns need to be wider for the new descriptions
    $ws.Columns.Item(16).ColumnWidth = 16.71
    $ws.Columns.Item(17).ColumnWidth = 16.71
}

# Reflect the editor's last selection (P5:Q5 on ScenarioB) while keeping
# ScenarioA as the active tab, matching the saved workbook view state.
$wsB = $wb.Worksheets.Item("ScenarioB")
$wsA = $wb.Worksheets.Item("ScenarioA")
$wsB.Select()
$wsB.Range("P5:Q5").Select()
$wsA.Select()
